$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 517, shifting existing rows 517-567 down to 518-568.
$ws.Rows.Item(517).Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(517, 1).Value = 7
$ws.Cells.Item(517, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(517, 3).Value = "Ñuble"
$ws.Cells.Item(517, 4).Value = 45166
$ws.Cells.Item(517, 5).Value = 16
$ws.Cells.Item(517, 6).Value = 100112002
$ws.Cells.Item(517, 7).Value = "Pimiento"
$ws.Cells.Item(517, 8).Value = "Zafiro verde"
$ws.Cells.Item(517, 9).Value = "Segunda"
$ws.Cells.Item(517, 10).Value = 50
$ws.Cells.Item(517, 11).Value = 12000
$ws.Cells.Item(517, 12).Value = 12000
$ws.Cells.Item(517, 13).Value = 12000
$ws.Cells.Item(517, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(517, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(517, 16).Value = 800
$ws.Cells.Item(517, 17).Value = 15
$ws.Cells.Item(517, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same date number formatting as the rest of column D.
$ws.Cells.Item(517, 4).NumberFormat = $ws.Cells.Item(518, 4).NumberFormat
